$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# C10 held the upper bound (18) for the R30 "Good Evening" rule; restore it
# to the value from the prior admin revision (1).
$ws.Range("C10").Value = 1
